$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws4 = $wb.Worksheets.Item(4)

# Sheet "展览" (sheet1) updates to column F ("想去人数")
$ws1.Range("F2").Value = 1358
$ws1.Range("F3").Value = 1229
$ws1.Range("F5").Value = 115
$ws1.Range("F7").Value = 674
$ws1.Range("F8").Value = 114
$ws1.Range("F11").Value = 2445
$ws1.Range("F12").Value = 1606
$ws1.Range("F13").Value = 1484
$ws1.Range("F14").Value = 312
$ws1.Range("F16").Value = 609
$ws1.Range("F17").Value = 785
$ws1.Range("F18").Value = 69
$ws1.Range("F19").Value = 310
$ws1.Range("F24").Value = 4991
$ws1.Range("F26").Value = 501
$ws1.Range("F27").Value = 77
$ws1.Range("F28").Value = 160
$ws1.Range("F29").Value = 137
$ws1.Range("F30").Value = 225
$ws1.Range("F31").Value = 127
$ws1.Range("F33").Value = 1039
$ws1.Range("F34").Value = 726
$ws1.Range("F36").Value = 52
$ws1.Range("F38").Value = 390
$ws1.Range("F39").Value = 1049
$ws1.Range("F42").Value = 170
$ws1.Range("F43").Value = 128
$ws1.Range("F44").Value = 24

# Sheet "全部类型" (sheet4) updates to column F ("想去人数")
$ws4.Range("F2").Value = 1358
$ws4.Range("F5").Value = 1230
$ws4.Range("F9").Value = 115
$ws4.Range("F11").Value = 674
$ws4.Range("F12").Value = 114
$ws4.Range("F17").Value = 2445
$ws4.Range("F18").Value = 1606
$ws4.Range("F19").Value = 1484
$ws4.Range("F20").Value = 312
$ws4.Range("F22").Value = 609
$ws4.Range("F24").Value = 785
$ws4.Range("F25").Value = 69
$ws4.Range("F26").Value = 310
$ws4.Range("F29").Value = 4991
$ws4.Range("F31").Value = 501
$ws4.Range("F32").Value = 77
$ws4.Range("F33").Value = 160
$ws4.Range("F34").Value = 137
$ws4.Range("F35").Value = 225
$ws4.Range("F36").Value = 128
$ws4.Range("F38").Value = 1039
$ws4.Range("F39").Value = 726
$ws4.Range("F40").Value = 52
$ws4.Range("F41").Value = 390
$ws4.Range("F42").Value = 1049
$ws4.Range("F44").Value = 170
$ws4.Range("F45").Value = 128
$ws4.Range("F46").Value = 24

$wb.Save()
